# Generate Report for handback
# This script updates the localization-status workbook to reflect that the
# two content files (cdd92b4f-... and f03929aa-...) have been handed back
# and are now in sync with en-US, for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the status shown for both locale columns
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Helper data describing the per-locale handback information:
#   sheetName        - the worksheet (table) name for the locale
#   handbackDateTime - the new "Latest Handback DateTime" for rows 2 & 3
# ---------------------------------------------------------------------
$locales = @(
    @{ Name = "zh-cn"; HandbackDateTime = "2016-02-19 06:19:03" },
    @{ Name = "de-de"; HandbackDateTime = "2016-02-19 06:19:21" }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Name)

    # Build a lookup of existing hyperlink addresses keyed by the range
    # address they are anchored to, so we can reuse them for the new
    # "Latest Target File" / "Latest Handback File" columns.
    $linkByAddress = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $rangeAddress = $hl.Range.Address()
        $linkByAddress[$rangeAddress] = $hl.Address
    }

    foreach ($row in 2, 3) {
        # Update the Status column (B)
        $ws.Range("B$row").Value = $newStatus

        $aCell = $ws.Range("A$row")
        $cCell = $ws.Range("C$row")

        $aText = $aCell.Value()
        $cText = $cCell.Value()

        $aUrl = $linkByAddress[$aCell.Address()]
        $cUrl = $linkByAddress[$cCell.Address()]

        # Latest Target File (E) mirrors the Source File Name (A) hyperlink,
        # and Latest Handback File (F) mirrors the Latest Handoff File (C)
        # hyperlink, since the file has been handed back unchanged / in sync.
        $ws.Hyperlinks.Add($ws.Range("E$row"), $aUrl, "", "", $aText)
        $ws.Hyperlinks.Add($ws.Range("F$row"), $cUrl, "", "", $cText)

        # Latest Handback DateTime (G)
        $ws.Range("G$row").Value = $locale.HandbackDateTime
    }
}

Write-Host "Handback report generated"
